# Fix the whitespace regex example in the proposal document.
#
# Before:  " \w+ (regular expression)"   (single run)
# After:   " \s"  +  "+ (regular expression)"   (split across two runs,
#          same character formatting: sz=20 / szCs=20)
#
# The original commit replaced the regex "\w" (word character) with "\s"
# (whitespace character) for the "whitespace" grammar rule. Because the
# edit happened in the middle of the existing run, Word (and this script)
# ends up with the replaced text living in its own run, separate from the
# trailing "+ (regular expression)" text.

$d = $word.ActiveDocument

# Locate " \w+ (regular expression)" (with its leading space) so we can
# work out exact character offsets for the sub-edits below.
$r = $d.Content
$found = $r.Find.Execute(" \w+ (regular expression)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $matchStart = $r.Start

    # Position of the "\w" (2 characters) right after the leading space.
    $wStart = $matchStart + 1
    $wRange = $d.Range($wStart, $wStart + 2)
    $wRange.Text = "\s"

    # The remaining text "+ (regular expression)" immediately follows.
    $restStart = $wStart + 2
    $restText = "+ (regular expression)"
    $restRange = $d.Range($restStart, $restStart + $restText.Length)

    # Toggling a character format on/off forces this trailing text into
    # its own run, separate from the " \s" run that precedes it, while
    # leaving the visible/effective formatting unchanged.
    $restRange.Bold = 1
    $restRange.Bold = 0
}
